$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.035.88"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "1.832.37"

$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'241.53"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").Value = "'0.6280"
$ws.Range("E6").Value = "  -4.54%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'45.59"
$ws.Range("E8").Value = "  +1.54%  "

$ws.Range("D9").Value = "'0.07570"
$ws.Range("E9").Value = "  +2.09%  "

$ws.Range("D10").Value = "'0.2919"
$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").Value = "'0.07647"
$ws.Range("E12").Value = "  -1.22%  "

$ws.Range("D13").Value = "1.828.12"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").Value = "'0.000009470"
$ws.Range("E17").Value = "  +9.94%  "

$ws.Range("D18").Value = "'5.987"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").Value = "28.956.19"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("D20").Value = "'225.39"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").Value = "'7.211"
$ws.Range("E23").Value = "  +1.43%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "'160.41"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").Value = "'8.434"
$ws.Range("E26").Value = "  -1.69%  "

$ws.Range("D27").Value = "'0.1366"
$ws.Range("E27").Value = "  -2.70%  "

$ws.Range("D28").Value = "'17.83"

$ws.Range("D29").Value = "'1.494"
$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("D30").Value = "'4.062"
$ws.Range("E30").Value = "  -0.91%  "

$ws.Range("D31").Value = "'4.032"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("E32").Value = "  +0.77%  "

$ws.Range("E33").Value = "  -1.10%  "

$ws.Range("D34").Value = "'1.848"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").Value = "'1.151"
$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("D36").Value = "'0.7330"

$ws.Range("D37").Value = "'2.594"
$ws.Range("E37").Value = "  -2.27%  "

$ws.Range("D38").Value = "1.271.47"
$ws.Range("E38").Value = "  -2.63%  "

$ws.Range("D39").Value = "'2.762"

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").Value = "'6.556"
$ws.Range("E41").Value = "  +7.86%  "

$ws.Range("D42").Value = "'0.8918"
$ws.Range("E42").Value = "  -2.57%  "

$ws.Range("D44").Value = "'101.81"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").Value = "1.977.39"
$ws.Range("E45").Value = "  -0.26%  "

$ws.Range("D46").Value = "'64.52"
$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("D49").Value = "'0.3981"
$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("D50").Value = "'8.839"
$ws.Range("E50").Value = "  +1.65%  "

$ws.Range("D51").Value = "'0.05756"
$ws.Range("E51").Value = "  -1.50%  "
